$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting D:K to E:L
$ws.Columns("D:D").Insert()

# Restore number formatting for the new column D by copying formats from column E
$ws.Range("E5:E102").Copy() | Out-Null
$ws.Range("D5:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate new column D with the new (most recent) fiscal-year figures
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 1207300
$ws.Range("D9").Value2 = 201100
$ws.Range("D10").Value2 = 1006200
$ws.Range("D12").Value2 = 48400
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = "NA"
$ws.Range("D17").Value2 = 963000
$ws.Range("D18").Value2 = 244300
$ws.Range("D20").Value2 = -66700
$ws.Range("D21").Value2 = 364800
$ws.Range("D22").Value2 = 0
$ws.Range("D23").Value2 = 177600
$ws.Range("D24").Value2 = 44800
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 132800
$ws.Range("D27").Value2 = 126800
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = 66700
$ws.Range("D33").Value2 = 126800
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 126800
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 209500
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value2 = 221600
$ws.Range("D44").Value2 = 0
$ws.Range("D45").Value2 = 29200
$ws.Range("D46").Value2 = 460300
$ws.Range("D47").Value2 = 83800
$ws.Range("D48").Value2 = 98800
$ws.Range("D49").Value2 = 1906800
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 11000
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 2560800
$ws.Range("D57").Value2 = 166500
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value2 = 140800
$ws.Range("D60").Value2 = 307300
$ws.Range("D61").Value2 = 1013100
$ws.Range("D62").Value2 = 204600
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 1525100
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 769600
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 1035700
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 126800
$ws.Range("D83").Value2 = 187200
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 401300
$ws.Range("D91").Value2 = -56400
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -406600
$ws.Range("D96").Value2 = -81700
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -131400
$ws.Range("D101").Value2 = -4800
$ws.Range("D102").Value2 = -141500
